$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D. Excel shifts the existing D:K quarterly
# data right to F:M, matching the added <xr:revisionPtr>-era resave captured by the diff
# (dimension grows from A5:K102 to A5:M102, spans grow from "x:11" to "x:13").
$ws.Range("D:E").Insert()

# The freshly inserted D:E columns come back blank/unstyled. Copy number format only
# (PasteSpecial xlPasteFormats = -4122) from the neighboring (now-shifted) original
# D:E data, which now lives in F:G, over to the new D:E cells -- scoped to the three
# contiguous data blocks so header-only rows (5, 6, 37, 79) are not touched.
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns (D, E) with the newest two reported quarters for every
# line item across the Income Statement, Balance Sheet and Cash Flow Statement blocks.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 79800
$ws.Range("E8").Value = 77400
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = -3700
$ws.Range("E15").Value = -3500
$ws.Range("D17").Value = 14400
$ws.Range("E17").Value = 12500
$ws.Range("D18").Value = 65400
$ws.Range("E18").Value = 64900
$ws.Range("D20").Value = -24400
$ws.Range("E20").Value = -23000
$ws.Range("D21").Value = 44800
$ws.Range("E21").Value = 45400
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 41100
$ws.Range("E23").Value = 41900
$ws.Range("D24").Value = 8300
$ws.Range("E24").Value = 9000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 32700
$ws.Range("E26").Value = 32900
$ws.Range("D27").Value = 32700
$ws.Range("E27").Value = 32900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 24400
$ws.Range("E32").Value = 23000
$ws.Range("D33").Value = 32700
$ws.Range("E33").Value = 32900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 32700
$ws.Range("E35").Value = 32900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 228400
$ws.Range("E41").Value = 185000
$ws.Range("D42").Value = 1195800
$ws.Range("E42").Value = 1607800
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 174400
$ws.Range("E48").Value = 170200
$ws.Range("D49").Value = 96200
$ws.Range("E49").Value = 97000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 7574300
$ws.Range("E54").Value = 7602400
$ws.Range("D57").Value = 37500
$ws.Range("E57").Value = 40300
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 26800
$ws.Range("E61").Value = 32000
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 6671500
$ws.Range("E66").Value = 6717600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 722600
$ws.Range("E72").Value = 707500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 902800
$ws.Range("E76").Value = 884800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 32700
$ws.Range("E81").Value = 32900
$ws.Range("D83").Value = 3700
$ws.Range("E83").Value = 3500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 31700
$ws.Range("E89").Value = 43700
$ws.Range("D91").Value = -7500
$ws.Range("E91").Value = -27800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -338700
$ws.Range("E94").Value = 36200
$ws.Range("D96").Value = -9800
$ws.Range("E96").Value = -6900
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -61000
$ws.Range("E100").Value = -55800
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -367900
$ws.Range("E102").Value = 24100
